$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(44803, 90, 24000, 24000, 24000, 1600)
    3 = @(44799, 80, 23000, 23000, 23000, 1533)
    4 = @(44831, 90, 25000, 25000, 25000, 1667)
    5 = @(44757, 80, 25000, 25000, 25000, 1667)
    6 = @(44789, 90, 24000, 24000, 24000, 1600)
    7 = @(44827, 90, 22000, 22000, 22000, 1467)
    8 = @(44761, 100, 23000, 25000, 24000, 1600)
    9 = @(44407, 90, 25000, 25000, 25000, 1667)
    10 = @(44817, 90, 23000, 23000, 23000, 1533)
    11 = @(44781, 70, 24000, 24000, 24000, 1600)
    12 = @(44418, 90, 25000, 25000, 25000, 1667)
    13 = @(44819, 70, 22000, 22000, 22000, 1467)
    14 = @(44838, 80, 22000, 22000, 22000, 1467)
    15 = @(44792, 120, 24000, 24000, 24000, 1600)
    16 = @(44740, 90, 25000, 25000, 25000, 1667)
    17 = @(44764, 90, 24000, 24000, 24000, 1600)
    18 = @(44775, 120, 24000, 24000, 24000, 1600)
    19 = @(44806, 70, 23000, 23000, 23000, 1533)
    20 = @(44750, 90, 25000, 25000, 25000, 1667)
    21 = @(44782, 120, 24000, 24000, 24000, 1600)
    22 = @(44771, 90, 25000, 25000, 25000, 1667)
    23 = @(44810, 110, 22000, 22000, 22000, 1467)
    24 = @(44365, 80, 25000, 25000, 25000, 1667)
    25 = @(44754, 90, 25000, 25000, 25000, 1667)
    26 = @(44778, 120, 24000, 24000, 24000, 1600)
    27 = @(44400, 80, 25000, 25000, 25000, 1667)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]   # D: Fecha
    $ws.Cells.Item($row, 10).Value = $vals[1]  # J: Volumen
    $ws.Cells.Item($row, 11).Value = $vals[2]  # K: Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals[3]  # L: Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals[4]  # M: Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals[5]  # P: Precio $/Kg
}
